# Add possible endings to the goals document
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Manieren om spel te eindigen" (ways to end the game) table ---
# Column A = Manier (way), Column B = Kamer (room code), Column C = Eindscenario (ending scenario number)

$data = @(
    @{ Row = 17; A = "Van terrein aflopen";                        B = "091100"; C = 1 },
    @{ Row = 18; A = "Nooduitgang bovenste verdieping gebruiken";   B = "070704"; C = 2 },
    @{ Row = 19; A = "Airco saboteren";                             B = "060704"; C = 3 },
    @{ Row = 20; A = "SSD-pc gebruiken";                            B = "060900"; C = 4 },
    @{ Row = 21; A = "IA-pc gebruiken";                             B = "050802"; C = 4 },
    @{ Row = 22; A = "Personeelsdossiers inzien";                   B = "070401"; C = 4 },
    @{ Row = 23; A = "Herrie maken in loods";                       B = "040600"; C = 4 },
    @{ Row = 24; A = "Sleutelbeen breken op serverruimtedeur";      B = "040600"; C = 5 },
    @{ Row = 25; A = "Kopieerapparaat misbruiken";                  B = "070802"; C = 6 }
)

# Column B in this table holds room codes such as "091100" that must stay text
# (leading zeros), so force a text number format before writing the values.
$ws.Range("B17:B28").NumberFormat = "@"

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# Left-align the whole table body (rows 17-28), matching the header above it.
$ws.Range("A17:C28").HorizontalAlignment = -4131

# --- Cosmetic/view updates that came along with the edit ---
# Column A needs to be a bit wider to fit the new, longer entries (target raw
# width 41.5703125 chars; ColumnWidth rounds to the nearest 1/6 of a
# character, so 40.6667 is the input that lands closest to that value).
$ws.Columns("A").ColumnWidth = 40.6667

# Selection moved to A20 while editing the new rows.
$ws.Range("A20").Select() | Out-Null
